$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6; existing rows 6-43 shift down to 7-44.
$ws.Range("A6:R6").EntireRow.Insert()

# Populate the new weekly entry in row 6.
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = "2022-05-04"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112030
$ws.Range("G6").Value = "Poroto granado"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 26000
$ws.Range("M6").Value = 25467
$ws.Range("N6").Value = "$/saco 25 kilos"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 1019
$ws.Range("Q6").Value = 25
$ws.Range("R6").Value = "Hortaliza"
